$d = $word.ActiveDocument

# Fix typo: {tracherName} -> {teacherName} in the signature/footer table.
# Locate the whole placeholder (braces included) and replace the Range's
# text directly so the run collapses into a single "{teacherName}" run
# (matching how Word coalesces runs/removes stale proofErr markers once
# the misspelling that triggered them is gone).
$rng = $d.Content
$find = $rng.Find
$find.ClearFormatting()
$find.Text = "{tracherName}"
$found = $find.Execute()
if ($found) {
    $rng.Text = "{teacherName}"
}
